# Generate Report for handoff
# - Adds a new "ready for handoff" file row (ffff30a78fcc-5af6-4ed1-85d5-ab7291c71512.md)
#   on every sheet, ahead of the ".localization-config" row.
# - Updates the existing file row's source filename/status from a failed
#   handoff to a successful one (3b651442-f62a-45d5-89af-f8da57062731.md /
#   "Ready for handoff"), and fills in the handoff file/datetime details.

$wb = $excel.ActiveWorkbook

$mdCommit     = "eee703022e0542320a8cd21e2d95497044864759"
$configCommit = "44dce391a3fff4d00e36c880da109017604ceace"

$file1 = "3b651442-f62a-45d5-89af-f8da57062731.md"
$file2 = "ffff30a78fcc-5af6-4ed1-85d5-ab7291c71512.md"
$configFile = ".localization-config"

$xlfZh = "3b651442-f62a-45d5-89af-f8da57062731.664a8cddf1addebc5b3a9066a8dd9c53bcb4a834.zh-cn.xlf"
$xlfDe = "3b651442-f62a-45d5-89af-f8da57062731.664a8cddf1addebc5b3a9066a8dd9c53bcb4a834.de-de.xlf"

$zeroDate = "0001-01-01 00:00:00"
$zhDate = "2016-02-16 15:27:05"
$deDate = "2016-02-16 15:27:18"

function Set-FileLink {
    param($ws, $cellRef, $fileName, $isConfig)

    $ws.Range($cellRef).Value = $fileName
    if ($isConfig) {
        $target = "https://github.com/OpenLocalizationTest/oltest/blob/$configCommit/$fileName"
    } else {
        $target = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$fileName"
    }
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, [Type]::Missing, [Type]::Missing, $fileName) | Out-Null
}

function Set-XlfLink {
    param($ws, $cellRef, $xlfName, $commitSha)

    $ws.Range($cellRef).Value = $xlfName
    $target = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSha/e2e/$xlfName"
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, [Type]::Missing, [Type]::Missing, $xlfName) | Out-Null
}

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Hyperlinks.Delete()

Set-FileLink $ov "A2" $file1 $false
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"

Set-FileLink $ov "A3" $file2 $false
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

Set-FileLink $ov "A4" $configFile $true
$ov.Range("B4").Value = "Not to be localized"
$ov.Range("C4").Value = "Not to be localized"

# ---------------------------------------------------------------------------
# zh-cn / de-de sheets share the same column layout:
# Source File Name | Status | Latest Handoff File | Latest Handoff Datetime |
# Latest Target File | Latest Handback File | Latest Handback DateTime |
# Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
function Update-LangSheet {
    param($ws, $xlfName, $xlfDate)

    $ws.Hyperlinks.Delete()

    # Row 2: the file that is now ready for handoff
    Set-FileLink $ws "A2" $file1 $false
    $ws.Range("B2").Value = "Ready for handoff"
    Set-XlfLink $ws "C2" $xlfName $mdCommit
    $ws.Range("D2").Value = $xlfDate
    $ws.Range("G2").Value = $zeroDate
    $ws.Range("H2").Value = "Include"

    # Row 3: new file row, also ready for handoff
    Set-FileLink $ws "A3" $file2 $false
    $ws.Range("B3").Value = "Ready for handoff"
    Set-XlfLink $ws "C3" $xlfName $mdCommit
    $ws.Range("D3").Value = $xlfDate
    $ws.Range("G3").Value = $zeroDate
    $ws.Range("H3").Value = "Include"

    # Row 4: the config file, not localized (shifted down from row 3)
    Set-FileLink $ws "A4" $configFile $true
    $ws.Range("B4").Value = "Not to be localized"
    $ws.Range("D4").Value = $zeroDate
    $ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("G4").Value = $zeroDate
    $ws.Range("H4").Value = "Ignored"
}

$zh = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $zh $xlfZh $zhDate

$de = $wb.Worksheets.Item("de-de")
Update-LangSheet $de $xlfDe $deDate
